# 23 dec 2023 update
# Fill in the Dec 17-22, 2023 loan-payment entries (columns J:L, rows 8-13)
# on the "MD10000.1-OCT" sheet, and move the active selection to reflect
# the newly-entered range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MD10000.1-OCT")

$dateFormat = "[$-409]d\-mmm\-yyyy;@"

# row -> (date serial, amount, count)
$entries = @(
    @{ Row = 8;  Date = 45277; Amount = 100; Count = 1 },
    @{ Row = 9;  Date = 45278; Amount = 100; Count = 1 },
    @{ Row = 10; Date = 45279; Amount = 100; Count = 1 },
    @{ Row = 11; Date = 45280; Amount = 100; Count = 1 },
    @{ Row = 12; Date = 45281; Amount = 100; Count = 1 },
    @{ Row = 13; Date = 45282; Amount = 100; Count = 1 }
)

foreach ($entry in $entries) {
    $r = $entry.Row

    $jCell = $ws.Cells.Item($r, 10)   # column J
    $jCell.Value = $entry.Date
    $jCell.NumberFormat = $dateFormat

    $kCell = $ws.Cells.Item($r, 11)   # column K
    $kCell.Value = $entry.Amount

    $lCell = $ws.Cells.Item($r, 12)   # column L
    $lCell.Value = $entry.Count
}

# Update the saved selection/active cell to match the newly-filled range.
$ws.Range("L7:L13").Select() | Out-Null
